# Mise à jour audit Excell
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 8 text content ---
$ws.Cells.Item(8, 3).Value2 = "Les div avec la classe ""keywords"" contient une succession de mots clés sans contexte pour tromper l'algo' Google, c'est une technique de blackhat.`nIl s'agit d'une technique de ""triche"" pour booster frauduleusement son référencement."
$ws.Cells.Item(8, 4).Value2 = "Essayer de tromper l'algorithme de Google représente un risque pour le référencement de la page car cela pourrait entrainer un malus.`nSuppression de ces divs là où elles apparaissent (header + footer donc)"

# --- Add new row 9 content (SEO et Accessibilité / Balises Alt) ---
$ws.Cells.Item(9, 1).Value2 = "SEO et Accessibilité"
$ws.Cells.Item(9, 2).Value2 = "Balises Alt "
$ws.Cells.Item(9, 3).Value2 = "La balises html Alt est certaines fois absente, mais quand elle est présente les mots cités ne décrivent pas le contenu de l'image. Cette balises est très utile pour les malvoyants mais aussi pour Google."
$ws.Cells.Item(9, 4).Value2 = "Ajout de la balise Alt sur les images qui ne l'ont pas et modification de leur contenu par des descriptions justes et concises."
$ws.Cells.Item(9, 5).Value2 = "X"
$ws.Cells.Item(9, 6).Value2 = "smartkeyword - Alt"

# Copy formatting from row 8 to row 9 so styles (fonts/wrap/etc.) match
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)  # xlPasteFormats

# Make sure the text values are still correct after paste (paste formats only touches formatting)
$ws.Cells.Item(9, 1).Value2 = "SEO et Accessibilité"
$ws.Cells.Item(9, 2).Value2 = "Balises Alt "
$ws.Cells.Item(9, 3).Value2 = "La balises html Alt est certaines fois absente, mais quand elle est présente les mots cités ne décrivent pas le contenu de l'image. Cette balises est très utile pour les malvoyants mais aussi pour Google."
$ws.Cells.Item(9, 4).Value2 = "Ajout de la balise Alt sur les images qui ne l'ont pas et modification de leur contenu par des descriptions justes et concises."
$ws.Cells.Item(9, 5).Value2 = "X"
$ws.Cells.Item(9, 6).Value2 = "smartkeyword - Alt"

$ws.Rows.Item(9).RowHeight = 36

# Add hyperlink on F9 pointing to a smartkeyword.io article about alt tags
$ws.Hyperlinks.Add($ws.Range("F9"), "https://smartkeyword.io/balise-alt-image-seo/")

# --- Update the sheet view (scroll position + selection) ---
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("F9").Select()
